$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H32").Value = 1795
$ws.Range("J32").Value = 2532.2
$ws.Range("L32").Value = 2532.2
$ws.Range("N32").Value = -3184.2
$ws.Range("H51").Value = 7096
$ws.Range("J51").Value = 7000
$ws.Range("L51").Value = 7000
$ws.Range("N51").Value = -7968
$ws.Range("H98").Value = 2323.7778
$ws.Range("I98").Value = 2323.7778
$ws.Range("K98").Value = 2323.7778
$ws.Range("M98").Value = -825.7777999999998
$ws.Range("H112").Value = 2901.7273
$ws.Range("J112").Value = 3020.8572
$ws.Range("L112").Value = 9062.571599999999
$ws.Range("N112").Value = -11278.5716
$ws.Range("H122").Value = 2323.7778
$ws.Range("I122").Value = 2323.7778
$ws.Range("K122").Value = 6971.3334
$ws.Range("M122").Value = -4521.3334
$ws.Range("H138").Value = 1919.8
$ws.Range("I138").Value = 1801.2142
$ws.Range("J138").Value = 2196.5
$ws.Range("K138").Value = 5403.642599999999
$ws.Range("L138").Value = 6589.5
$ws.Range("M138").Value = -263.6425999999992
$ws.Range("N138").Value = -16869.5

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 3781.8215
$ws.Range("I32").Value = 2365.0715
$ws.Range("K32").Value = 2365.0715
$ws.Range("M32").Value = -2078.0715
$ws.Range("H45").Value = 1630.7894
$ws.Range("I45").Value = 956.4286
$ws.Range("K45").Value = 956.4286
$ws.Range("M45").Value = -579.4286
$ws.Range("H74").Value = 6306.3335
$ws.Range("I74").Value = 4459.5
$ws.Range("K74").Value = 4459.5
$ws.Range("M74").Value = -3585.5
$ws.Range("H77").Value = 6306.3335
$ws.Range("I77").Value = 4459.5
$ws.Range("K77").Value = 22297.5
$ws.Range("M77").Value = -17929.5
$ws.Range("H132").Value = 1853.2
$ws.Range("I132").Value = 1349.1428
$ws.Range("J132").Value = 4499.5
$ws.Range("K132").Value = 4047.4284
$ws.Range("L132").Value = 13498.5
$ws.Range("M132").Value = -1517.4284
$ws.Range("N132").Value = -18558.5

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H44").Value = 22000
$ws.Range("J44").Value = 22000
$ws.Range("L44").Value = 22000
$ws.Range("N44").Value = -22994
$ws.Range("H105").Value = 2493.5625
$ws.Range("I105").Value = 2493.5625
$ws.Range("K105").Value = 2493.5625
$ws.Range("M105").Value = -746.5625
$ws.Range("H107").Value = 2247.4
$ws.Range("I107").Value = 1699.6666
$ws.Range("J107").Value = 2482.1428
$ws.Range("K107").Value = 1699.6666
$ws.Range("L107").Value = 2482.1428
$ws.Range("M107").Value = 220.3334
$ws.Range("N107").Value = -6322.1428

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 3685.6667
$ws.Range("I31").Value = 1378.3
$ws.Range("J31").Value = 5783.273
$ws.Range("K31").Value = 1378.3
$ws.Range("L31").Value = 5783.273
$ws.Range("M31").Value = -1083.3
$ws.Range("N31").Value = -6373.273
$ws.Range("H34").Value = 3685.6667
$ws.Range("I34").Value = 1378.3
$ws.Range("J34").Value = 5783.273
$ws.Range("K34").Value = 1378.3
$ws.Range("L34").Value = 5783.273
$ws.Range("M34").Value = -1176.3
$ws.Range("N34").Value = -6187.273
$ws.Range("H63").Value = 35000
$ws.Range("J63").Value = 35000
$ws.Range("L63").Value = 35000
$ws.Range("N63").Value = -36372
$ws.Range("H66").Value = 35000
$ws.Range("J66").Value = 35000
$ws.Range("L66").Value = 105000
$ws.Range("N66").Value = -111864
$ws.Range("H107").Value = 499.26315
$ws.Range("J107").Value = 606.8570999999999
$ws.Range("L107").Value = 606.8570999999999
$ws.Range("N107").Value = -4446.8571
$ws.Range("H134").Value = 1053.6154
$ws.Range("I134").Value = 1016.7273
$ws.Range("J134").Value = 1256.5
$ws.Range("K134").Value = 3050.1819
$ws.Range("L134").Value = 3769.5
$ws.Range("M134").Value = -515.1819
$ws.Range("N134").Value = -8839.5

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 563.6429000000001
$ws.Range("I5").Value = 498.83334
$ws.Range("K5").Value = 1496.50002
$ws.Range("M5").Value = -1384.50002
$ws.Range("H14").Value = 115.72727
$ws.Range("I14").Value = 115.72727
$ws.Range("K14").Value = 347.18181
$ws.Range("M14").Value = -174.18181
$ws.Range("H56").Value = 17191.166
$ws.Range("I56").Value = 17191.166
$ws.Range("K56").Value = 17191.166
$ws.Range("M56").Value = -16661.166
$ws.Range("H113").Value = 7480.2
$ws.Range("I113").Value = 33768
$ws.Range("J113").Value = 908.25
$ws.Range("K113").Value = 101304
$ws.Range("L113").Value = 2724.75
$ws.Range("M113").Value = -99134
$ws.Range("N113").Value = -7064.75
$ws.Range("H131").Value = 13535781
$ws.Range("J131").Value = 28259.172
$ws.Range("L131").Value = 84777.516
$ws.Range("N131").Value = -94857.516
$ws.Range("H132").Value = 1544.625
$ws.Range("J132").Value = 1519.5
$ws.Range("L132").Value = 13675.5
$ws.Range("N132").Value = -18735.5
$ws.Range("H135").Value = 563.6429000000001
$ws.Range("I135").Value = 498.83334
$ws.Range("K135").Value = 4489.50006
$ws.Range("M135").Value = -1954.50006

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H11").Value = 5848147
$ws.Range("I11").Value = 8752220
$ws.Range("J11").Value = 40001.6
$ws.Range("K11").Value = 8752220
$ws.Range("L11").Value = 40001.6
$ws.Range("M11").Value = -8752081
$ws.Range("N11").Value = -40279.6
$ws.Range("H132").Value = 2825.389
$ws.Range("I132").Value = 2510.1667
$ws.Range("J132").Value = 4401.5
$ws.Range("K132").Value = 7530.500100000001
$ws.Range("L132").Value = 13204.5
$ws.Range("M132").Value = -5000.500100000001
$ws.Range("N132").Value = -18264.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 5992.0713
$ws.Range("I7").Value = 3000.1667
$ws.Range("K7").Value = 3000.1667
$ws.Range("M7").Value = -2888.1667
$ws.Range("H82").Value = 3552.9092
$ws.Range("I82").Value = 1867.3334
$ws.Range("K82").Value = 1867.3334
$ws.Range("M82").Value = -1506.3334
$ws.Range("H85").Value = 3552.9092
$ws.Range("I85").Value = 1867.3334
$ws.Range("K85").Value = 1867.3334
$ws.Range("M85").Value = -619.3334
$ws.Range("H126").Value = 5992.0713
$ws.Range("I126").Value = 3000.1667
$ws.Range("K126").Value = 9000.500100000001
$ws.Range("M126").Value = -6530.500100000001
$ws.Range("H136").Value = 4840.4
$ws.Range("I136").Value = 2134.6667
$ws.Range("K136").Value = 6404.000100000001
$ws.Range("M136").Value = -3854.000100000001

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H126").Value = 4108.0312
$ws.Range("I126").Value = 3293.2727
$ws.Range("J126").Value = 5900.5
$ws.Range("K126").Value = 9879.8181
$ws.Range("L126").Value = 17701.5
$ws.Range("M126").Value = -7409.8181
$ws.Range("N126").Value = -22641.5
$ws.Range("H132").Value = 2815.4
$ws.Range("I132").Value = 878.25
$ws.Range("J132").Value = 6689.7
$ws.Range("K132").Value = 2634.75
$ws.Range("L132").Value = 20069.1
$ws.Range("M132").Value = -104.75
$ws.Range("N132").Value = -25129.1
